# Update the 5x5 grid of two-digit ÷ one-digit division problems.
# The worksheet table has 20 rows (every 4th row, i.e. rows 1,5,9,13,17,
# holds the actual answer text; the rows in-between are blank spacer rows).
# Each of those 5 "answer" rows has 5 cells (columns 1-5).
#
# We address every cell by its row/column position (rather than doing a
# global Find/Replace by old text) because some of the new values equal
# old values used elsewhere in the grid - e.g. "45÷2=22, 1" is both a
# target value and, later, a source value; a blind Replace-All would
# clobber itself. Writing straight to Cell.Range.Text avoids that and
# keeps each run's existing formatting (rFonts/sz) untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "12÷5=2, 2"
$t.Rows.Item(1).Cells.Item(2).Range.Text = "70÷8=8, 6"
$t.Rows.Item(1).Cells.Item(3).Range.Text = "45÷2=22, 1"
$t.Rows.Item(1).Cells.Item(4).Range.Text = "18÷9=2, 0"
$t.Rows.Item(1).Cells.Item(5).Range.Text = "94÷9=10, 4"

$t.Rows.Item(5).Cells.Item(1).Range.Text = "95÷8=11, 7"
$t.Rows.Item(5).Cells.Item(2).Range.Text = "18÷6=3, 0"
$t.Rows.Item(5).Cells.Item(3).Range.Text = "97÷7=13, 6"
$t.Rows.Item(5).Cells.Item(4).Range.Text = "11÷3=3, 2"
$t.Rows.Item(5).Cells.Item(5).Range.Text = "88÷2=44, 0"

$t.Rows.Item(9).Cells.Item(1).Range.Text = "44÷7=6, 2"
$t.Rows.Item(9).Cells.Item(2).Range.Text = "31÷3=10, 1"
$t.Rows.Item(9).Cells.Item(3).Range.Text = "21÷7=3, 0"
$t.Rows.Item(9).Cells.Item(4).Range.Text = "34÷3=11, 1"
$t.Rows.Item(9).Cells.Item(5).Range.Text = "16÷2=8, 0"

$t.Rows.Item(13).Cells.Item(1).Range.Text = "41÷8=5, 1"
$t.Rows.Item(13).Cells.Item(2).Range.Text = "33÷7=4, 5"
$t.Rows.Item(13).Cells.Item(3).Range.Text = "11÷7=1, 4"
$t.Rows.Item(13).Cells.Item(4).Range.Text = "44÷6=7, 2"
$t.Rows.Item(13).Cells.Item(5).Range.Text = "46÷2=23, 0"

$t.Rows.Item(17).Cells.Item(1).Range.Text = "75÷9=8, 3"
$t.Rows.Item(17).Cells.Item(2).Range.Text = "67÷2=33, 1"
$t.Rows.Item(17).Cells.Item(3).Range.Text = "24÷8=3, 0"
$t.Rows.Item(17).Cells.Item(4).Range.Text = "81÷6=13, 3"
$t.Rows.Item(17).Cells.Item(5).Range.Text = "77÷4=19, 1"

Write-Output "Updated 25 cells across 5 rows."
